$d = $word.ActiveDocument

# 1. Replace the placeholder text (including the trailing-space run) in the
#    first paragraph with the new placeholder, collapsing the two runs into one.
$null = $d.Content.Find.Execute("**ID__AFFARS_5312_topic_3__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5312_102__ID**", 2)

# 2. Update the first paragraph's formatting: add a paragraph border (5pt
#    space on every side, no visible line) and change the left indent from
#    120 twips (6pt) to 225 twips (11.25pt).
$p1 = $d.Paragraphs(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25
